# Horarios actualizados Línea 141 - 728
# Applies the scraped-schedule refresh (new scrape time 05:27:50) to all
# three worksheets: LP1912, LP1912-215, 6203-6173.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LP1912")

$ws.Range("A2").Value = "Última actualización: 05:27:50"
$ws.Range("A3").Value = "Total filas: 42"

# Updates to existing rows (Hora_Scrap / Minutos refreshed by the new
# scrape; a couple of rows also swapped Hora_Scrap/Linea/Minutos values).
$ws.Cells.Item(23,1).Value = "05:27:50"
$ws.Cells.Item(23,4).Value = 7

$ws.Cells.Item(24,1).Value = "03:42:43"
$ws.Cells.Item(24,3).Value = "14_ABASTO"
$ws.Cells.Item(24,4).Value = 113

$ws.Cells.Item(25,1).Value = "04:17:03"
$ws.Cells.Item(25,3).Value = "215B_EL PATO"
$ws.Cells.Item(25,4).Value = 78

$ws.Cells.Item(27,1).Value = "05:27:50"
$ws.Cells.Item(27,4).Value = 19

$ws.Cells.Item(28,1).Value = "05:27:50"
$ws.Cells.Item(28,4).Value = 27

$ws.Cells.Item(29,1).Value = "05:27:50"
$ws.Cells.Item(29,4).Value = 37

$ws.Cells.Item(31,1).Value = "05:27:50"
$ws.Cells.Item(31,4).Value = 44

$ws.Cells.Item(33,1).Value = "05:27:50"
$ws.Cells.Item(33,4).Value = 47

$ws.Cells.Item(34,1).Value = "05:27:50"
$ws.Cells.Item(34,4).Value = 54

$ws.Cells.Item(35,1).Value = "05:27:50"
$ws.Cells.Item(35,4).Value = 60

$ws.Cells.Item(36,1).Value = "05:27:50"
$ws.Cells.Item(36,4).Value = 62

$ws.Cells.Item(37,1).Value = "05:27:50"
$ws.Cells.Item(37,4).Value = 64

$ws.Cells.Item(38,1).Value = "05:27:50"
$ws.Cells.Item(38,4).Value = 77

$ws.Cells.Item(39,1).Value = "05:27:50"
$ws.Cells.Item(39,4).Value = 79

# New arrivals appended by the refresh.
$newRows1 = @(
    @("05:27:50","06:59","14_ABASTO",92,"LP1912"),
    @("05:27:50","07:04","23_HERNANDEZ",97,"LP1912"),
    @("05:27:50","07:05","15_ABASTO",98,"LP1912"),
    @("05:27:50","07:06","225_GOMEZ",99,"LP1912"),
    @("05:27:50","07:11","215A_EL PATO",104,"LP1912"),
    @("05:27:50","07:15","11_ETCHEVERRY",108,"LP1912"),
    @("05:27:50","07:21","26_HERNANDEZ",114,"LP1912"),
    @("05:27:50","07:23","10_OLMOS",116,"LP1912")
)
$r = 40
foreach ($row in $newRows1) {
    $ws.Cells.Item($r,1).Value = $row[0]
    $ws.Cells.Item($r,2).Value = $row[1]
    $ws.Cells.Item($r,3).Value = $row[2]
    $ws.Cells.Item($r,4).Value = $row[3]
    $ws.Cells.Item($r,5).Value = $row[4]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LP1912-215")

$ws.Range("A2").Value = "Última actualización: 05:27:50"
$ws.Range("A3").Value = "Total filas: 15"

$ws.Cells.Item(15,1).Value = "05:27:50"
$ws.Cells.Item(15,4).Value = 7

$ws.Cells.Item(17,1).Value = "05:27:50"
$ws.Cells.Item(17,4).Value = 44

$ws.Cells.Item(19,1).Value = "05:27:50"
$ws.Cells.Item(19,4).Value = 79

$ws.Cells.Item(20,1).Value = "05:27:50"
$ws.Cells.Item(20,2).Value = "07:11"
$ws.Cells.Item(20,3).Value = "215A_EL PATO"
$ws.Cells.Item(20,4).Value = 104
$ws.Cells.Item(20,5).Value = "LP1912"

# ---------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("6203-6173")

$ws.Range("A2").Value = "Última actualización: 05:27:50"
$ws.Range("A3").Value = "Total filas: 8"

$ws.Cells.Item(7,1).Value = "05:27:50"
$ws.Cells.Item(7,4).Value = 16

$ws.Cells.Item(9,1).Value = "05:27:50"
$ws.Cells.Item(9,4).Value = 41

$ws.Cells.Item(11,1).Value = "05:27:50"
$ws.Cells.Item(11,4).Value = 65

$ws.Cells.Item(13,1).Value = "05:27:50"
$ws.Cells.Item(13,2).Value = "06:59"
$ws.Cells.Item(13,3).Value = "215B_LP-P MOR-1 Y 57"
$ws.Cells.Item(13,4).Value = 92
$ws.Cells.Item(13,5).Value = "L6173"
